$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old "category" labels (column D) to their renamed counterparts.
# The underlying model/std/mean data (columns A-C) and the model names
# (column E) are untouched -- only the category text is swapped for a new
# (more verb-like) label.
$categoryMap = @{
    "statement pos"     = "make statement pos"
    "cooperation pos"   = "cooperate pos"
    "retreat pos"       = "yield pos"
    "investigation pos" = "investigate pos"
    "demand pos"        = "demand pos"
    "dissaproval pos"   = "disapprove pos"
    "rejection pos"     = "reject pos"
    "threat pos"        = "threaten pos"
    "protest pos"       = "protest pos"
    "force pos"         = "exhibit force pos"
    "relation pos"      = "reduce relations pos"
    "coercion pos"      = "coerce pos"
    "assault pos"       = "assault pos"
    "fight pos"         = "fight pos"
    "hybrid attack pos" = "mass violence pos"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    if ($categoryMap.ContainsKey($current)) {
        $cell.Value = $categoryMap[$current]
    }
}

# Resize column D to fit the new (longer) category labels, matching the
# author's re-run of AutoFit after the rename ("reduce relations pos" is now
# the longest label, so the bestFit column grows wider).
$ws.Columns.Item(4).ColumnWidth = 18.65

# Reset the view: scroll back to the top and move the active selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("I10").Select() | Out-Null
